# Remove Argentina (and its now-orphaned "India" entry that was dropped in
# the same edit) from the population list, then refresh the sheet's cached
# sort-range state and restore the active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Argentina" row (exact-match lookup, re-verified before deleting
# so a stray partial match can never take out the wrong country).
$argentinaCell = $ws.Columns.Item(1).Find("Argentina")
if (($argentinaCell -ne $null) -and ($argentinaCell.Text -eq "Argentina")) {
    $ws.Rows.Item($argentinaCell.Row).EntireRow.Delete() | Out-Null
}

# Delete the "India" row (same defensive exact-match check; guards against
# "India" being treated as a partial match of "Indonesia").
$indiaCell = $ws.Columns.Item(1).Find("India")
if (($indiaCell -ne $null) -and ($indiaCell.Text -eq "India")) {
    $ws.Rows.Item($indiaCell.Row).EntireRow.Delete() | Out-Null
}

# The table (A2:D82, 81 countries) is already alphabetically sorted by
# Country once the two rows are gone, but re-apply the sort over the sheet's
# previously-remembered sort range (A2:C80 / A2:A80) so that cached range
# shrinks by the same two rows we removed, matching the saved sort state
# (A2:C78 / A2:A78).
$sortObj = $ws.Sort
$sortObj.SortFields.Clear() | Out-Null
$sortObj.SortFields.Add($ws.Range("A2:A78")) | Out-Null
$sortObj.SetRange($ws.Range("A2:C78")) | Out-Null
$sortObj.Header = 2
$sortObj.Apply() | Out-Null

# Leave the active cell on A2, matching the saved selection state.
$ws.Range("A2").Select() | Out-Null
